$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to text so values like "58.536.08" or "537.15"
# are stored verbatim instead of being reinterpreted as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "58.536.08"
$ws.Range("E2").Value = "  -0.22%  "
$ws.Range("D3").Value = "2.302.69"
$ws.Range("E3").Value = "  +0.47%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "537.15"
$ws.Range("E5").Value = "  -1.76%  "
$ws.Range("D6").Value = "132.26"
$ws.Range("E6").Value = "  +1.34%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  +2.51%  "
$ws.Range("D9").Value = "2.300.36"
$ws.Range("E9").Value = "  +0.44%  "
$ws.Range("D10").Value = "0.101"
$ws.Range("E10").Value = "  -1.06%  "
$ws.Range("E11").Value = "  -0.62%  "
$ws.Range("E12").Value = "  +0.87%  "
$ws.Range("D13").Value = "0.335"
$ws.Range("E13").Value = "  -0.11%  "
$ws.Range("D14").Value = "23.83"
$ws.Range("E14").Value = "  +0.29%  "
$ws.Range("D15").Value = "2.713.44"
$ws.Range("E15").Value = "  +0.50%  "
$ws.Range("D16").Value = "58.475.92"
$ws.Range("E16").Value = "  -0.28%  "
$ws.Range("E17").Value = "  -0.24%  "
$ws.Range("D18").Value = "2.297.47"
$ws.Range("E18").Value = "  +1.55%  "
$ws.Range("D19").Value = "10.59"
$ws.Range("E19").Value = "  -0.53%  "
$ws.Range("D20").Value = "4.21"
$ws.Range("E20").Value = "  -1.97%  "
$ws.Range("D21").Value = "316.31"
$ws.Range("E21").Value = "  +0.49%  "
$ws.Range("D22").Value = "6.64"
$ws.Range("E22").Value = "  +2.65%  "
$ws.Range("E23").Value = "  +0.25%  "
$ws.Range("D24").Value = "63.12"
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("E25").Value = "  -1.45%  "
$ws.Range("E26").Value = "  +0.15%  "
$ws.Range("E27").Value = "  -1.45%  "
$ws.Range("E28").Value = "  -1.06%  "
$ws.Range("D29").Value = "170.95"
$ws.Range("E29").Value = "  +0.79%  "
$ws.Range("E30").Value = "  -1.97%  "
$ws.Range("D31").Value = "0.0₃0726"
$ws.Range("E31").Value = "  +0.05%  "
$ws.Range("E32").Value = "  +1.94%  "
$ws.Range("E33").Value = "  +0.68%  "
$ws.Range("D34").Value = "0.385"
$ws.Range("E34").Value = "  +0.34%  "
$ws.Range("D36").Value = "17.89"
$ws.Range("E36").Value = "  +0.58%  "
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("D38").Value = "1.25"
$ws.Range("E38").Value = "  -0.09%  "
$ws.Range("D39").Value = "4.01"
$ws.Range("E39").Value = "  +1.39%  "
$ws.Range("E40").Value = "  -0.17%  "
$ws.Range("D41").Value = "291.66"
$ws.Range("E41").Value = "  -2.62%  "
$ws.Range("D42").Value = "140.98"
$ws.Range("E42").Value = "  +0.52%  "
$ws.Range("E43").Value = "  +0.37%  "
$ws.Range("D44").Value = "0.0953"
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").Value = "0.0497"
$ws.Range("E45").Value = "  -0.47%  "
$ws.Range("D46").Value = "0.557"
$ws.Range("E46").Value = "  +0.15%  "
$ws.Range("D47").Value = "18.36"
$ws.Range("E47").Value = "  -1.00%  "
$ws.Range("E48").Value = "  -1.91%  "
$ws.Range("E49").Value = "  -0.43%  "
$ws.Range("E51").Value = "  +0.85%  "

# Restore the default (unstyled) look for column D now that the text is committed,
# matching the workbook's original formatting (no explicit style on these cells).
$ws.Range("D2:D51").Style = "Normal"
